# Experience curve + formula.xlsx - workbook restructuring
# - Reorder sheets: Experience, d2 mean mlvl, Monsters, Damage Types, Attack Type, Enemy_Modifiers (new)
# - Add tab colors to every sheet
# - Add a "Type" column header to the Monsters table
# - Update selections / active sheet
# - Add a brand-new "Enemy_Modifiers" sheet with data
#
# NOTE: worksheet object handles in this runtime stay bound to their
# positional slot, not to sheet identity - so after any Move/reorder we
# must re-fetch the worksheet we want **by name** before touching it again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Reorder sheets so "d2 mean mlvl" sits between "Experience" and "Monsters"
# ---------------------------------------------------------------------------
$wsMonsters = $wb.Worksheets.Item("Monsters")
$wsD2 = $wb.Worksheets.Item("d2 mean mlvl")
$wsMonsters.Move($null, $wsD2)

# Re-fetch by name now that positions have shifted.
$wsD2 = $wb.Worksheets.Item("d2 mean mlvl")
$wsMonsters = $wb.Worksheets.Item("Monsters")

# ---------------------------------------------------------------------------
# 2. Tab colours for the existing sheets
# ---------------------------------------------------------------------------
$wsExperience = $wb.Worksheets.Item("Experience")
$wsExperience.Tab.Color = 15773696   # FF00B0F0

$wsD2.Tab.Color = 15773696           # FF00B0F0
$wsMonsters.Tab.Color = 5287936      # FF00B050

$wsDamage = $wb.Worksheets.Item("Damage Types")
$wsDamage.Tab.Color = 5287936        # FF00B050

$wsAttack = $wb.Worksheets.Item("Attack Type")
$wsAttack.Tab.Color = 5287936        # FF00B050

# ---------------------------------------------------------------------------
# 3. Add the "Type" header to the Monsters sheet, next to "Name"
# ---------------------------------------------------------------------------
$wsMonsters = $wb.Worksheets.Item("Monsters")
$wsMonsters.Activate()
$wsMonsters.Range("C1").Value = "Type"
$wsMonsters.Range("A1:O1").Select()

# ---------------------------------------------------------------------------
# 4. Update the selection on the Experience sheet
# ---------------------------------------------------------------------------
$wsExperience = $wb.Worksheets.Item("Experience")
$wsExperience.Activate()
$wsExperience.Range("H36").Select()

# ---------------------------------------------------------------------------
# 5. Add the new "Enemy_Modifiers" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "Enemy_Modifiers"
$wsNew = $wb.Worksheets.Item("Enemy_Modifiers")
$wsNew.Tab.Color = 65535             # FFFFFF00

$wsNew.Range("A1").Value = "id"
$wsNew.Range("B1").Value = "Name"
$wsNew.Range("C1").Value = "Type"
$wsNew.Range("D1").Value = "damage_type"
$wsNew.Range("E1").Value = "Modifier_Type"
$wsNew.Range("F1").Value = "Modifier"
$wsNew.Range("G1").Value = "Rarity (chance of effect occurring)"

$wsNew.Range("B2").Value = "Chilled"
$wsNew.Range("C2").Value = "damage"
$wsNew.Range("D2").Value = "elemental-Ice"
$wsNew.Range("E2").Value = "Additive"
$wsNew.Range("F2").Value = 20
$wsNew.Range("G2").Value = 25

$wsNew.Range("B3").Value = "Frozen"
$wsNew.Range("C3").Value = "damage"
$wsNew.Range("D3").Value = "elemental-Ice"
$wsNew.Range("E3").Value = "Multiplicative"
$wsNew.Range("F3").Value = 1.5
$wsNew.Range("G3").Value = 10

$wsNew.Range("B4").Value = "Flaming"
$wsNew.Range("C4").Value = "damage"
$wsNew.Range("D4").Value = "elemental-fire"
$wsNew.Range("E4").Value = "Additive"
$wsNew.Range("F4").Value = 20
$wsNew.Range("G4").Value = 25

$wsNew.Range("B5").Value = "Burning"
$wsNew.Range("C5").Value = "damage"
$wsNew.Range("D5").Value = "elemental-fire"
$wsNew.Range("E5").Value = "Multiplicative"
$wsNew.Range("F5").Value = 1.5
$wsNew.Range("G5").Value = 10

$wsNew.Range("B6").Value = "Rotten"
$wsNew.Range("C6").Value = "damage"
$wsNew.Range("D6").Value = "elemental-poison"
$wsNew.Range("E6").Value = "Additive"
$wsNew.Range("F6").Value = 20
$wsNew.Range("G6").Value = 25

$wsNew.Range("B7").Value = "Necrotic"
$wsNew.Range("C7").Value = "damage"
$wsNew.Range("D7").Value = "elemental-poison"
$wsNew.Range("E7").Value = "Multiplicative"
$wsNew.Range("F7").Value = 1.5
$wsNew.Range("G7").Value = 10

$wsNew.Range("B8").Value = "Lucrative"
$wsNew.Range("C8").Value = "chanceOfItem"
$wsNew.Range("D8").Value = "null"
$wsNew.Range("E8").Value = "Multiplicative"
$wsNew.Range("F8").Value = 1.5
$wsNew.Range("G8").Value = 5

$wsNew.Range("B9").Value = "Rich"
$wsNew.Range("C9").Value = "chanceOfGold"
$wsNew.Range("D9").Value = "null"
$wsNew.Range("E9").Value = "Multiplicative"
$wsNew.Range("F9").Value = 10
$wsNew.Range("G9").Value = 5

# Column widths roughly matching the authored layout (best-fit on text columns)
$wsNew.Columns.Item(3).AutoFit()
$wsNew.Columns.Item(4).AutoFit()
$wsNew.Columns.Item(5).ColumnWidth = $wsNew.Columns.Item(4).ColumnWidth
$wsNew.Columns.Item(6).AutoFit()

$wsNew.Range("G10").Select()
